$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$urls = @(
    "https://www.mercadolivre.com.br/apple-iphone-14-128-gb-amarelo-distribuidor-autorizado/p/MLB1024156591?pdp_filters=category:MLB1055#searchVariation=MLB1024156591&position=2&search_layout=stack&type=product&tracking_id=c9fc395a-5c68-4423-a1b3-0cf8f16d183f",
    "https://www.mercadolivre.com.br/apple-iphone-13-128-gb-estelar-distribuidor-autorizado/p/MLB1018500855?pdp_filters=category:MLB1055#searchVariation=MLB1018500855&position=3&search_layout=stack&type=product&tracking_id=c9fc395a-5c68-4423-a1b3-0cf8f16d183f",
    "https://www.mercadolivre.com.br/apple-iphone-14-128-gb-meia-noite-distribuidor-autorizado/p/MLB1019615360?pdp_filters=category:MLB1055#searchVariation=MLB1019615360&position=4&search_layout=stack&type=product&tracking_id=c9fc395a-5c68-4423-a1b3-0cf8f16d183f",
    "https://www.mercadolivre.com.br/apple-iphone-13-128-gb-meia-noite-distribuidor-autorizado/p/MLB1018500844?pdp_filters=category:MLB1055#searchVariation=MLB1018500844&position=5&search_layout=stack&type=product&tracking_id=c9fc395a-5c68-4423-a1b3-0cf8f16d183f",
    "https://www.mercadolivre.com.br/apple-iphone-14-128-gb-estelar-distribuidor-autorizado/p/MLB1019615378?pdp_filters=category:MLB1055#searchVariation=MLB1019615378&position=6&search_layout=stack&type=product&tracking_id=c9fc395a-5c68-4423-a1b3-0cf8f16d183f"
)

$startRow = 218
for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $urls[$i]
}
